$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- New row 134: "How to Think Like a Roman Emperor" ---
# Copy an existing full data row first so number formats / styles carry over.
$ws.Range("A131:I131").Copy($ws.Range("A134:I134"))
$ws.Range("A134").Value = "How to Think Like a Roman Emperor"
$ws.Range("B134").Value = "Donald Robertson"
$ws.Range("C134").Value = 44124
$ws.Range("D134").Value = 44126
$ws.Range("E134").Value = "stoic;marcus aurelius;philosophy;history"
$ws.Range("F134").Value = "Audio"
$ws.Range("G134").Value = "8 Hours 36 Mins"
$ws.Range("H134").Value = 3
$ws.Range("I134").Value = $true

# --- New row 135: "168 Hours" ---
$ws.Range("A131:I131").Copy($ws.Range("A135:I135"))
$ws.Range("A135").Value = "168 Hours"
$ws.Range("B135").Value = "Laura Vanderkam"
$ws.Range("C135").Value = 44126
$ws.Range("D135").Value = 44128
$ws.Range("E135").Value = "productivity;time management;achievement;personal improvement"
$ws.Range("F135").Value = "Audio"
$ws.Range("G135").Value = "7 Hours 48 Mins"
$ws.Range("H135").Value = 4
$ws.Range("I135").Value = $true

# Row 133 ("Covid-19") previously had no value in the "First Time Reading?"
# column; the edit also backfills it as TRUE.
$ws.Range("I133").Value = $true

# Move the active selection to the next empty row below the new data,
# matching where Excel leaves the cursor after entering the new rows.
$ws.Range("A136").Select()
